# maj planning + script mur
# Update the "Planning" tracker: mark several tasks as done (100%) or
# partially done, roll back one task to 0%, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Progress cells that move from 0% (red, style s=4) to 100% (green, style s=6) ---
# Copy the fill/number-format from an existing "100%" cell (B29) so the
# written style reuses the workbook's existing "done" style instead of
# creating a new one, then set the value to 1 (100%).
$doneSrc = $ws.Range("B29")
$doneTargets = @("B9", "B10", "B11", "B31")
foreach ($ref in $doneTargets) {
    $doneSrc.Copy()
    $ws.Range($ref).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($ref).Value = 1
}

# --- Progress cells that move from 0% (red) to a partial percentage (orange, style s=7) ---
$partialSrc = $ws.Range("B30")

$partialSrc.Copy()
$ws.Range("B33").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B33").Value = 0.8

$partialSrc.Copy()
$ws.Range("B37").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B37").Value = 0.2

$excel.CutCopyMode = $false

# --- Task that regresses from 100% back to 0% (style stays the same, s=6) ---
$ws.Range("B46").Value = 0

# --- Move the active selection to where work resumed ---
$ws.Range("E9").Select() | Out-Null
